# Applies the cell-value edits described by the OOXML diff for cryptos.xlsx.
# All changed cells are plain text (inline strings) in the source workbook -
# some "Price" values look like plain numbers (e.g. "259.65", "1.00"), and if
# assigned via .Value directly, Excel/COM auto-converts them to numeric cells
# (losing formatting like trailing zeros, e.g. "1.00" -> 1). To preserve the
# original text semantics we force text format before the write, then clear
# the (now unneeded) formatting override so no stray style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '37.329.13'
Set-TextValue $ws.Range("E2") '  -0.13%  '
# Row 3
Set-TextValue $ws.Range("D3") '2.009.71'
Set-TextValue $ws.Range("E3") '  -1.27%  '
# Row 5
Set-TextValue $ws.Range("D5") '259.65'
Set-TextValue $ws.Range("E5") '  +4.43%  '
# Row 6
Set-TextValue $ws.Range("D6") '0.617'
Set-TextValue $ws.Range("E6") '  -2.05%  '
# Row 7
Set-TextValue $ws.Range("E7") '  -0.02%  '
# Row 8
Set-TextValue $ws.Range("D8") '56.53'
Set-TextValue $ws.Range("E8") '  -7.29%  '
# Row 9
Set-TextValue $ws.Range("D9") '0.382'
Set-TextValue $ws.Range("E9") '  -3.60%  '
# Row 10
Set-TextValue $ws.Range("D10") '0.0772'
Set-TextValue $ws.Range("E10") '  -5.14%  '
# Row 11
Set-TextValue $ws.Range("E11") '  -3.28%  '
# Row 12
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D12") '2.307.08'
Set-TextValue $ws.Range("E12") '  -1.25%  '
# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D13") '14.25'
Set-TextValue $ws.Range("E13") '  -7.41%  '
# Row 14
Set-TextValue $ws.Range("D14") '21.64'
Set-TextValue $ws.Range("E14") '  -3.92%  '
# Row 15
Set-TextValue $ws.Range("D15") '0.793'
Set-TextValue $ws.Range("E15") '  -8.20%  '
# Row 16
Set-TextValue $ws.Range("D16") '5.19'
Set-TextValue $ws.Range("E16") '  -6.50%  '
# Row 17
Set-TextValue $ws.Range("D17") '2.003.47'
Set-TextValue $ws.Range("E17") '  -1.47%  '
# Row 18
Set-TextValue $ws.Range("D18") '37.247.91'
Set-TextValue $ws.Range("E18") '  -0.19%  '
# Row 19
Set-TextValue $ws.Range("D19") '69.98'
Set-TextValue $ws.Range("E19") '  -1.25%  '
# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0835'
Set-TextValue $ws.Range("E20") '  -4.16%  '
# Row 21
Set-TextValue $ws.Range("D21") '232.29'
Set-TextValue $ws.Range("E21") '  +0.30%  '
# Row 22
Set-TextValue $ws.Range("D22") '5.10'
Set-TextValue $ws.Range("E22") '  -3.41%  '
# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D23") '1.00'
Set-TextValue $ws.Range("E23") '  -0.08%  '
# Row 24
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D24") '2.59'
Set-TextValue $ws.Range("E24") '  +0.66%  '
# Row 25
Set-TextValue $ws.Range("E25") '  -0.39%  '
# Row 26
Set-TextValue $ws.Range("D26") '164.55'
Set-TextValue $ws.Range("E26") '  +0.34%  '
# Row 27
Set-TextValue $ws.Range("D27") '8.94'
Set-TextValue $ws.Range("E27") '  -6.26%  '
# Row 28
Set-TextValue $ws.Range("D28") '19.54'
Set-TextValue $ws.Range("E28") '  -1.72%  '
# Row 29
Set-TextValue $ws.Range("E29") '  -6.13%  '
# Row 30
Set-TextValue $ws.Range("D30") '1.31'
# Row 31
Set-TextValue $ws.Range("E31") '  -2.22%  '
# Row 32
Set-TextValue $ws.Range("D32") '4.58'
Set-TextValue $ws.Range("E32") '  -5.77%  '
# Row 33
Set-TextValue $ws.Range("D33") '0.0638'
Set-TextValue $ws.Range("E33") '  -5.11%  '
# Row 34
Set-TextValue $ws.Range("D34") '4.46'
Set-TextValue $ws.Range("E34") '  -1.90%  '
# Row 35
Set-TextValue $ws.Range("D35") '2.36'
Set-TextValue $ws.Range("E35") '  -6.33%  '
# Row 36
Set-TextValue $ws.Range("E36") '  +0.30%  '
# Row 37
Set-TextValue $ws.Range("E37") '  +0.13%  '
# Row 38
Set-TextValue $ws.Range("E38") '  -8.48%  '
# Row 39
Set-TextValue $ws.Range("D39") '5.44'
Set-TextValue $ws.Range("E39") '  -0.76%  '
# Row 40
Set-TextValue $ws.Range("E40") '  +1.70%  '
# Row 41
Set-TextValue $ws.Range("D41") '1.18'
Set-TextValue $ws.Range("E41") '  -1.28%  '
# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D42") '0.0211'
Set-TextValue $ws.Range("E42") '  -1.92%  '
# Row 43
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D43") '0.0923'
Set-TextValue $ws.Range("E43") '  -6.32%  '
# Row 44
Set-TextValue $ws.Range("D44") '1.433.43'
Set-TextValue $ws.Range("E44") '  +3.02%  '
# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D45") '89.37'
Set-TextValue $ws.Range("E45") '  -4.08%  '
# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D46") '15.68'
Set-TextValue $ws.Range("E46") '  -8.84%  '
# Row 47
Set-TextValue $ws.Range("D47") '1.02'
Set-TextValue $ws.Range("E47") '  -4.15%  '
# Row 48
Set-TextValue $ws.Range("E48") '  +2.32%  '
# Row 49
Set-TextValue $ws.Range("D49") '7.00'
Set-TextValue $ws.Range("E49") '  -6.81%  '
# Row 50
Set-TextValue $ws.Range("D50") '2.200.44'
Set-TextValue $ws.Range("E50") '  -1.19%  '
# Row 51
Set-TextValue $ws.Range("E51") '  -11.36%  '
